# Update recomputed TPM-based NATMI ligand-receptor metrics (Lipc-Lrp1).
# Only numeric value cells change; labels in columns A-D, K, L are untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.7303676666666666
$ws.Range("H2").Value = 2.191103
$ws.Range("I2").Value = 0.4916846149668853
$ws.Range("J2").Value = 0.4916846149668851
$ws.Range("M2").Value = 3.795192333333334
$ws.Range("N2").Value = 11.385577
$ws.Range("O2").Value = 0.01044213755712683
$ws.Range("P2").Value = 0.01044213755712683
$ws.Range("Q2").Value = 2.771885769047889
$ws.Range("R2").Value = 24.946971921431
$ws.Range("S2").Value = 0.005134238384207158
$ws.Range("T2").Value = 0.005134238384207158

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.7303676666666666
$ws.Range("H3").Value = 2.191103
$ws.Range("I3").Value = 0.4916846149668853
$ws.Range("J3").Value = 0.4916846149668851
$ws.Range("N3").Value = 730.1291960000001
$ws.Range("O3").Value = 0.6696287328350964
$ws.Range("P3").Value = 0.6696287328350964
$ws.Range("Q3").Value = 177.7542524159098
$ws.Range("R3").Value = 1599.788271743188
$ws.Range("S3").Value = 0.3292461456747877
$ws.Range("T3").Value = 0.3292461456747875

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.7303676666666666
$ws.Range("H4").Value = 2.191103
$ws.Range("I4").Value = 0.4916846149668853
$ws.Range("J4").Value = 0.4916846149668851
$ws.Range("M4").Value = 29.801371
$ws.Range("N4").Value = 89.404113
$ws.Range("O4").Value = 0.08199584844219236
$ws.Range("P4").Value = 0.08199584844219235
$ws.Range("Q4").Value = 21.76595780073767
$ws.Range("R4").Value = 195.893620206639
$ws.Range("S4").Value = 0.04031609717018243
$ws.Range("T4").Value = 0.04031609717018241

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.7303676666666666
$ws.Range("H5").Value = 2.191103
$ws.Range("I5").Value = 0.4916846149668853
$ws.Range("J5").Value = 0.4916846149668851
$ws.Range("M5").Value = 86.47679266666667
$ws.Range("N5").Value = 259.430378
$ws.Range("O5").Value = 0.2379332811655844
$ws.Range("P5").Value = 0.2379332811655844
$ws.Range("Q5").Value = 63.15985328077044
$ws.Range("R5").Value = 568.438679526934
$ws.Range("S5").Value = 0.116988133737708
$ws.Range("T5").Value = 0.116988133737708

# Row 6
$ws.Range("G6").Value = 0.2633103333333333
$ws.Range("H6").Value = 0.7899309999999999
$ws.Range("I6").Value = 0.1772609136062552
$ws.Range("J6").Value = 0.1772609136062552
$ws.Range("M6").Value = 3.795192333333334
$ws.Range("N6").Value = 11.385577
$ws.Range("O6").Value = 0.01044213755712683
$ws.Range("P6").Value = 0.01044213755712683
$ws.Range("Q6").Value = 0.9993133583541112
$ws.Range("R6").Value = 8.993820225187001
$ws.Range("S6").Value = 0.001850982843378493
$ws.Range("T6").Value = 0.001850982843378492

# Row 7
$ws.Range("G7").Value = 0.2633103333333333
$ws.Range("H7").Value = 0.7899309999999999
$ws.Range("I7").Value = 0.1772609136062552
$ws.Range("J7").Value = 0.1772609136062552
$ws.Range("N7").Value = 730.1291960000001
$ws.Range("O7").Value = 0.6696287328350964
$ws.Range("P7").Value = 0.6696287328350964
$ws.Range("Q7").Value = 64.08352065838622
$ws.Range("S7").Value = 0.1186990009593482
$ws.Range("T7").Value = 0.1186990009593482

# Row 8
$ws.Range("G8").Value = 0.2633103333333333
$ws.Range("H8").Value = 0.7899309999999999
$ws.Range("I8").Value = 0.1772609136062552
$ws.Range("J8").Value = 0.1772609136062552
$ws.Range("M8").Value = 29.801371
$ws.Range("N8").Value = 89.404113
$ws.Range("O8").Value = 0.08199584844219236
$ws.Range("P8").Value = 0.08199584844219235
$ws.Range("Q8").Value = 7.847008931800333
$ws.Range("R8").Value = 70.62308038620299
$ws.Range("S8").Value = 0.01453465900678306
$ws.Range("T8").Value = 0.01453465900678305

# Row 9
$ws.Range("G9").Value = 0.2633103333333333
$ws.Range("H9").Value = 0.7899309999999999
$ws.Range("I9").Value = 0.1772609136062552
$ws.Range("J9").Value = 0.1772609136062552
$ws.Range("M9").Value = 86.47679266666667
$ws.Range("N9").Value = 259.430378
$ws.Range("O9").Value = 0.2379332811655844
$ws.Range("P9").Value = 0.2379332811655844
$ws.Range("Q9").Value = 22.77023310265755
$ws.Range("R9").Value = 204.932097923918
$ws.Range("S9").Value = 0.04217627079674549
$ws.Range("T9").Value = 0.04217627079674548

# Row 10
$ws.Range("G10").Value = 0.331774
$ws.Range("H10").Value = 0.995322
$ws.Range("I10").Value = 0.223350757284377
$ws.Range("J10").Value = 0.2233507572843769
$ws.Range("M10").Value = 3.795192333333334
$ws.Range("N10").Value = 11.385577
$ws.Range("O10").Value = 0.01044213755712683
$ws.Range("P10").Value = 0.01044213755712683
$ws.Range("Q10").Value = 1.259146141199334
$ws.Range("R10").Value = 11.332315270794
$ws.Range("S10").Value = 0.002332259331051913
$ws.Range("T10").Value = 0.002332259331051912

# Row 11
$ws.Range("G11").Value = 0.331774
$ws.Range("H11").Value = 0.995322
$ws.Range("I11").Value = 0.223350757284377
$ws.Range("J11").Value = 0.2233507572843769
$ws.Range("N11").Value = 730.1291960000001
$ws.Range("O11").Value = 0.6696287328350964
$ws.Range("P11").Value = 0.6696287328350964
$ws.Range("Q11").Value = 80.74596129123468
$ws.Range("R11").Value = 726.7136516211121
$ws.Range("S11").Value = 0.1495620845780965
$ws.Range("T11").Value = 0.1495620845780965

# Row 12
$ws.Range("G12").Value = 0.331774
$ws.Range("H12").Value = 0.995322
$ws.Range("I12").Value = 0.223350757284377
$ws.Range("J12").Value = 0.2233507572843769
$ws.Range("M12").Value = 29.801371
$ws.Range("N12").Value = 89.404113
$ws.Range("O12").Value = 0.08199584844219236
$ws.Range("P12").Value = 0.08199584844219235
$ws.Range("Q12").Value = 9.887320062154
$ws.Range("R12").Value = 88.985880559386
$ws.Range("S12").Value = 0.01831383484373867
$ws.Range("T12").Value = 0.01831383484373866

# Row 13
$ws.Range("G13").Value = 0.331774
$ws.Range("H13").Value = 0.995322
$ws.Range("I13").Value = 0.223350757284377
$ws.Range("J13").Value = 0.2233507572843769
$ws.Range("M13").Value = 86.47679266666667
$ws.Range("N13").Value = 259.430378
$ws.Range("O13").Value = 0.2379332811655844
$ws.Range("P13").Value = 0.2379332811655844
$ws.Range("Q13").Value = 28.69075141019067
$ws.Range("R13").Value = 258.216762691716
$ws.Range("S13").Value = 0.05314257853148986
$ws.Range("T13").Value = 0.05314257853148985

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.1599873333333333
$ws.Range("H14").Value = 0.4799620000000001
$ws.Range("I14").Value = 0.1077037141424827
$ws.Range("J14").Value = 0.1077037141424827
$ws.Range("M14").Value = 3.795192333333334
$ws.Range("N14").Value = 11.385577
$ws.Range("O14").Value = 0.01044213755712683
$ws.Range("P14").Value = 0.01044213755712683
$ws.Range("Q14").Value = 0.6071827008971112
$ws.Range("R14").Value = 5.464644308074002
$ws.Range("S14").Value = 0.001124656998489271
$ws.Range("T14").Value = 0.001124656998489271

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.1599873333333333
$ws.Range("H15").Value = 0.4799620000000001
$ws.Range("I15").Value = 0.1077037141424827
$ws.Range("J15").Value = 0.1077037141424827
$ws.Range("N15").Value = 730.1291960000001
$ws.Range("O15").Value = 0.6696287328350964
$ws.Range("P15").Value = 0.6696287328350964
$ws.Range("Q15").Value = 38.93714101895023
$ws.Range("R15").Value = 350.4342691705521
$ws.Range("S15").Value = 0.07212150162286413
$ws.Range("T15").Value = 0.07212150162286411

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.1599873333333333
$ws.Range("H16").Value = 0.4799620000000001
$ws.Range("I16").Value = 0.1077037141424827
$ws.Range("J16").Value = 0.1077037141424827
$ws.Range("M16").Value = 29.801371
$ws.Range("N16").Value = 89.404113
$ws.Range("O16").Value = 0.08199584844219236
$ws.Range("P16").Value = 0.08199584844219235
$ws.Range("Q16").Value = 4.767841875967333
$ws.Range("R16").Value = 42.910576883706
$ws.Range("S16").Value = 0.00883125742148822
$ws.Range("T16").Value = 0.008831257421488216

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.1599873333333333
$ws.Range("H17").Value = 0.4799620000000001
$ws.Range("I17").Value = 0.1077037141424827
$ws.Range("J17").Value = 0.1077037141424827
$ws.Range("M17").Value = 86.47679266666667
$ws.Range("N17").Value = 259.430378
$ws.Range("O17").Value = 0.2379332811655844
$ws.Range("P17").Value = 0.2379332811655844
$ws.Range("Q17").Value = 13.83519145395956
$ws.Range("R17").Value = 124.516723085636
$ws.Range("S17").Value = 0.02562629809964105
$ws.Range("T17").Value = 0.02562629809964105

Write-Output "Applied 203 cell updates"
